# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the data rows (2-43) of the single data
# table on the sheet: every row's full record (columns A:R) moves to a new
# row position. Row 18 is unchanged. We read every source row into memory
# first (so reads never see already-overwritten data), then write each
# record out to its destination row per the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (1-based worksheet row numbers)
$map = @{
    2  = 12
    3  = 41
    4  = 26
    5  = 14
    6  = 42
    7  = 27
    8  = 25
    9  = 40
    10 = 34
    11 = 24
    12 = 10
    13 = 38
    14 = 13
    15 = 2
    16 = 7
    17 = 22
    18 = 18
    19 = 30
    20 = 21
    21 = 32
    22 = 31
    23 = 36
    24 = 43
    25 = 15
    26 = 16
    27 = 11
    28 = 17
    29 = 8
    30 = 6
    31 = 19
    32 = 29
    33 = 28
    34 = 35
    35 = 37
    36 = 9
    37 = 39
    38 = 23
    39 = 3
    40 = 20
    41 = 5
    42 = 33
    43 = 4
}

# Snapshot every source row (A:R) before any writes happen.
$snapshot = @{}
foreach ($r in $map.Keys) {
    $srcRow = $map[$r]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rng = $ws.Range("A" + $srcRow + ":R" + $srcRow)
        $snapshot[$srcRow] = $rng.Value2
    }
}

# Write each destination row from its recorded source snapshot.
foreach ($r in $map.Keys) {
    $srcRow = $map[$r]
    $destRng = $ws.Range("A" + $r + ":R" + $r)
    $destRng.Value = $snapshot[$srcRow]
}
